$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.403.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "'3.141.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'603.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "'143.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'3.137.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "'5.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.51%  "
$ws.Range("D14").Value = "'35.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'3.665.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "'64.432.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "'3.139.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'479.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "'14.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'0.711"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("D23").Value = "'7.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'85.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "'13.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "'8.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  +8.81%  "
$ws.Range("D30").Value = "'2.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("D31").Value = "'0.114"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Value = "'26.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'2.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").Value = "'1.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "'0.0₃0775"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.99%  "
$ws.Range("D37").Value = "'5.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("D38").Value = "'52.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'3.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.40%  "
$ws.Range("D40").Value = "'444.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").Value = "'0.0395"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").Value = "'0.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "'8.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'2.852.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").Value = "'0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "'2.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D49").Value = "'26.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "'120.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.23%  "
